# UPDATE technology portfolios for Norway
# Updates p_ites_c_inv (E2) base value on sheet "2025" from 5200 to 4600,
# which cascades via formulas into the other year sheets (2030-2050),
# and updates p_ites_c_charge_discharge (H2) from 5.5 to 5 on every sheet.

$wb = $excel.ActiveWorkbook

# Update the base investment cost on the 2025 sheet; downstream sheets
# reference this cell via formulas ('2025'!E2*(1-0.25*x)) so they will
# recalculate automatically.
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("E2").Value = 4600
$ws2025.Range("H2").Value = 5

# Update the charge/discharge cost on every other year sheet as well.
$sheetNames = @("2030", "2035", "2040", "2045", "2050")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("H2").Value = 5
}

$excel.Calculate()
